# Script produced by the data-refresh scraper run on 07-01-2024.
# - Re-orders a few match rows that were scraped out of their intended
#   chronological position (the two fixtures of a same match-day were
#   written in swapped order) by exchanging all match-specific columns
#   (F..V, except the opening-odds timestamps K/O/S which are identical
#   for both fixtures of the same match-day and stay put).
# - Appends two newly scraped fixtures at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-MatchRows {
    param($ws, [int]$r1, [int]$r2)

    $cols = @("F","G","H","I","J","L","M","N","P","Q","R","T","U","V")

    $row1 = @{}
    $row2 = @{}
    foreach ($c in $cols) {
        $row1[$c] = $ws.Range("$c$r1").Value2
        $row2[$c] = $ws.Range("$c$r2").Value2
    }
    foreach ($c in $cols) {
        $ws.Range("$c$r1").Value = $row2[$c]
        $ws.Range("$c$r2").Value = $row1[$c]
    }
}

# --- Swap the two fixtures that were stored in the wrong order ---
Swap-MatchRows $ws 15 16
Swap-MatchRows $ws 74 75
Swap-MatchRows $ws 96 97

# --- Append the two newly scraped fixtures (rows 114 & 115) ---
$ws.Range("A113:V113").Copy($ws.Range("A114:V114")) | Out-Null
$ws.Range("A113:V113").Copy($ws.Range("A115:V115")) | Out-Null

$ws.Range("A114").Value = 113
$ws.Range("B114").Value = "morocco"
$ws.Range("C114").Value = "botola-pro"
$ws.Range("D114").Value = "2023-2024"
$ws.Range("E114").Value = 45298.66666666666
$ws.Range("F114").Value = "Hassania Agadir"
$ws.Range("G114").Value = 3
$ws.Range("H114").Value = "Union Touarga"
$ws.Range("I114").Value = 2
$ws.Range("J114").Value = 3.03
$ws.Range("K114").Value = "04/01/2024 23:42"
$ws.Range("L114").Value = 2.19
$ws.Range("M114").Value = "07/01/2024 15:59"
$ws.Range("N114").Value = 2.96
$ws.Range("O114").Value = "04/01/2024 23:42"
$ws.Range("P114").Value = 3.09
$ws.Range("Q114").Value = "07/01/2024 15:59"
$ws.Range("R114").Value = 2.29
$ws.Range("S114").Value = "04/01/2024 23:42"
$ws.Range("T114").Value = 3.43
$ws.Range("U114").Value = "07/01/2024 15:59"
$ws.Range("V114").Value = "https://www.betexplorer.com/football/morocco/botola-pro/hassania-agadir-union-touarga/CKZCU3AI/"

$ws.Range("A115").Value = 114
$ws.Range("B115").Value = "morocco"
$ws.Range("C115").Value = "botola-pro"
$ws.Range("D115").Value = "2023-2024"
$ws.Range("E115").Value = 45298.66666666666
$ws.Range("F115").Value = "Renaissance Zemamra"
$ws.Range("G115").Value = 2
$ws.Range("H115").Value = "Berkane"
$ws.Range("I115").Value = 1
$ws.Range("J115").Value = 4.13
$ws.Range("K115").Value = "04/01/2024 22:12"
$ws.Range("L115").Value = 5.39
$ws.Range("M115").Value = "07/01/2024 15:57"
$ws.Range("N115").Value = 2.96
$ws.Range("O115").Value = "04/01/2024 22:12"
$ws.Range("P115").Value = 3.18
$ws.Range("Q115").Value = "07/01/2024 15:57"
$ws.Range("R115").Value = 1.89
$ws.Range("S115").Value = "04/01/2024 22:12"
$ws.Range("T115").Value = 1.75
$ws.Range("U115").Value = "07/01/2024 15:57"
$ws.Range("V115").Value = "https://www.betexplorer.com/football/morocco/botola-pro/renaissance-zemamra-berkane/jBYGTNPO/"

Write-Output "rows updated: swapped 15/16, 74/75, 96/97; appended 114, 115"
